$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new "MuLan" reference entry (title + body paragraphs) right
#    before the existing "Improving information retrieval by semantic
#    embedding" paragraph. We anchor on the END of the previous paragraph
#    (the BERT architecture paragraph) so the insertion becomes clean,
#    independent paragraphs rather than merging into the next paragraph.
# ---------------------------------------------------------------------------
$anchorFind = $d.Content.Find
$anchorFind.Execute("linearly projected embedding space.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorRange = $anchorFind.Parent
$insertPoint = $d.Range($anchorRange.End, $anchorRange.End)

$newParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:pStyle w:val="2"/>
        <w:bidi w:val="0"/>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>MULAN: A JOINT EMBEDDING OF MUSIC AUDIO AND NATURAL LANGUAGE</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Music tagging and content-based retrieval systems have traditionally been constructed using pre-defined ontologies covering a rigid set of music attributes or text queries. This paper presents MuLan: a first attempt at a new generation of acoustic models that link music audio directly to unconstrained natural language music descriptions. MuLan takes the form of a two-tower, joint audio-text embedding model trained using 44 million music recordings (370K hours) and weakly-associated, free-form text annotations. Through its compatibility with a wide range of music genres and text styles (including conventional music tags), the resulting audio-text representation subsumes existing ontologies while graduating to true zero-shot functionalities. We demonstrate the versatility of the MuLan embeddings with a range of experiments including transfer learning, zero-shot music tagging, language understanding in the music domain, and cross-modal retrieval applications.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Classifiers are generally trained to label examples with predefined and fixed class inventories, which are often manually specified as a structured ontology indicating interclass relationships. Empowered by recent advances in neural language modeling and their demonstrated transfer learning competence, researchers have begun exploring less restrictive natural language interfaces to access</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>the categorical information underlying raw content signals. The majority of this work has been in the visual and audio event domain, where a recent series of studies have demonstrated the utility of jointly embedding media content with natural language captions [1–5].  These joint embeddings have demonstrated strong capabilities in a range of applications, including transfer learning, cross-modal retrieval, automatic captioning, and zero-shot classification.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>The success of these efforts strongly depends on large-scale training resources and hefty neural network architectures that are flexible enough to model the complex, non-monotonic relationship between language and other modalities. In particular, the visual domain has greatly benefited from the availability of large amounts of captioned images available across the web [1]. However, in</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>the general environmental audio domain, such large-scale audio-caption pairs are less readily available and related efforts have relied on small captioned datasets [6, 7]. Critically, these datasets do not span the diversity of sound-descriptive language and their success in the more difficult zero-shot setting has been lacking [3, 8, 9].</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>our strategy is to assemble a collection of textual annotations extracted from metadata, comments, and playlist data and map them to a training set of over 44 million internet music videos. As was the case with image-text model training in [1], our text data only truly refers to the musical content in a fraction of cases. Therefore, we also explore text pre-filtering using a text classifier separately trained to identify music descriptions.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Music text joint embedding models</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Content-based music information retrieval requires linking the rich semantics expressible to free-form text with both broad and finegrained musical properties. One approach is to consider a large number of text label classes and try to ground the semantics in music with a multi-label classification task …..</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>our two-tower parallel encoder approach results in a joint embedding space that provides a natural language interface to arbitrary music audio. This opens up downstream opportunities for cross-modal retrieval, zero-shot tagging, and language understanding.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">Evaluation: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="400" w:firstLineChars="200"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Music Retrieval from Text Queries</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Given a music search collection and a text query, MuLan provides the ability to retrieve the music clips that are closest to the query in the embedding space. This evaluation is relevant to music retrieval applications, where content features can offer finer-grained and more complete similarity information when compared with metadata-based methods [41].</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Results and Discussion</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="400" w:firstLineChars="200"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:u w:val="single"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Music Retrieval from Text Queries</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Even though we start with a BERT checkpoint pretrained with massive language resources, training MuLan with only AudioSet clips and label annotations provides very limited ability to ground indomain natural language to music. Such limited crossmodal supervision does not generalize to the rich semantics that appear in the playlist titles and descriptions, which are more in line with the complex queries that are presented to real-world music search engines. We observe significant gain after including the large-scale short-form tags mined from the internet, which helps the model learn to ground more fine-grained music concepts. There i</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>s additional gain when including comments and playlist data, where the complete sentences are helpful for grounding the more complex queries, including multi-term queries (e.g.‘instrumental action movie soundtrack’), compositional queries (e.g. ‘classical music with middle eastern influence’), and even queries with negation (e.g. ‘hard rock without vocals’). Again, we find that training is surprisingly robust to annotation noise, achieving similar performance using unfiltered training text.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="2"/>
        <w:bidi w:val="0"/>
        <w:spacing w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Contrastive audio-language learning for music</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>As one of the most intuitive interfaces known to humans, natural language has the potential to mediate many tasks that involve human-computer interaction, especially in application-focused fields like Music Information Retrieval.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Our approach consists of a dual-encoder architecture that learns the alignment between pairs of music audio and descriptive sentences, producing multimodal embeddings that can be used for text-to-audio and audio-to-text retrieval out-of-the-box.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from its old location (end of the
#    Document-Retrieval-and-Question-Answering-with-BERT entry) to the new
#    edit location inside the freshly typed MuLan text, splitting
#    "There i" | "s additional gain ...".
# ---------------------------------------------------------------------------
$goBackFind = $d.Content.Find
$goBackFind.Execute("There i", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRange = $goBackFind.Parent
$goBackPoint = $d.Range($goBackRange.End, $goBackRange.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# ---------------------------------------------------------------------------
# 3. Update the "Normal (Web)" style default font to Times New Roman /
#    SimSun (matching the rest of the document's default run fonts).
# ---------------------------------------------------------------------------
$webStyle = $d.Styles("Normal (Web)")
$webStyle.Font.Name = "Times New Roman"
$webStyle.Font.NameFarEast = "SimSun"
$webStyle.Font.NameBi = "Times New Roman"

Write-Output "done"
